# Fix apply_defaults for biquad modules
#
# Row 23 (GPIO 21, "edge detector" section / trigger related): fill in the
# n_bits / description / module columns that were previously blank, and
# flag it with a note in the new column I.
#
# Row 35 (first biquad_0 parameter row): fill in the n_bits / description /
# module columns that were previously blank (pre-amp, 4 bits,
# coarse_gain_and_limiter), matching the pattern used by sibling rows
# below it (e.g. row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order these new strings are first written in determines their
# position in the shared-string table, so write them in the same order as
# the target workbook (pre-amp, new, trig mode).
$ws.Range("F35").Value = "pre-amp"
$ws.Range("I23").Value = "new"
$ws.Range("F23").Value = "trig mode"

$ws.Range("E23").Value = 1
$ws.Range("G23").Value = "mux_2x1"

$ws.Range("E35").Value = 4
$ws.Range("G35").Value = "coarse_gain_and_limiter"

# I35 stays empty but picks up the row's standard formatting (same style
# as the rest of row 35, e.g. H35).
$ws.Range("I35").Style = $ws.Range("H35").Style

# Matches the author's last selection in the saved file.
[void]$ws.Range("I23").Select()
